# Apply cached market-price / profit recalculations to the Anima Profits workbook.
# Values below are static (non-formula) numeric snapshots; update in place per sheet.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 1034.3334
$ws.Range("I29").Value = 403
$ws.Range("J29").Value = 1350
$ws.Range("K29").Value = 1209
$ws.Range("L29").Value = 4050
$ws.Range("M29").Value = -928
$ws.Range("N29").Value = -4612
$ws.Range("H38").Value = 272.44446
$ws.Range("I38").Value = 194
$ws.Range("K38").Value = 582
$ws.Range("M38").Value = -210
$ws.Range("H40").Value = 2020.8334
$ws.Range("I40").Value = 1822.8572
$ws.Range("K40").Value = 1822.8572
$ws.Range("M40").Value = -1647.8572
$ws.Range("H58").Value = 1129.2
$ws.Range("I58").Value = 911.5
$ws.Range("J58").Value = 2000
$ws.Range("K58").Value = 2734.5
$ws.Range("L58").Value = 6000
$ws.Range("M58").Value = -2584.5
$ws.Range("N58").Value = -6300
$ws.Range("H64").Value = 3722.111
$ws.Range("I64").Value = 3812.375
$ws.Range("J64").Value = 3000
$ws.Range("K64").Value = 3812.375
$ws.Range("L64").Value = 3000
$ws.Range("M64").Value = -3564.375
$ws.Range("N64").Value = -3496
$ws.Range("H67").Value = 3722.111
$ws.Range("I67").Value = 3812.375
$ws.Range("J67").Value = 3000
$ws.Range("K67").Value = 3812.375
$ws.Range("L67").Value = 3000
$ws.Range("M67").Value = -2954.375
$ws.Range("N67").Value = -4716
$ws.Range("H92").Value = 55559080
$ws.Range("I92").Value = 66670096
$ws.Range("K92").Value = 66670096
$ws.Range("M92").Value = -66668848
$ws.Range("H100").Value = 1814.238
$ws.Range("I100").Value = 1628.1428
$ws.Range("J100").Value = 2186.4285
$ws.Range("K100").Value = 1628.1428
$ws.Range("L100").Value = 2186.4285
$ws.Range("M100").Value = -1087.1428
$ws.Range("N100").Value = -3268.4285
$ws.Range("H115").Value = 1966.9
$ws.Range("I115").Value = 1780.7142
$ws.Range("K115").Value = 5342.142599999999
$ws.Range("M115").Value = -3775.142599999999
$ws.Range("H132").Value = 3907.6
$ws.Range("I132").Value = 3661.524
$ws.Range("J132").Value = 5199.5
$ws.Range("K132").Value = 10984.572
$ws.Range("L132").Value = 15598.5
$ws.Range("M132").Value = -8454.572
$ws.Range("N132").Value = -20658.5
$ws.Range("H138").Value = 2274.8103
$ws.Range("I138").Value = 1792.92
$ws.Range("J138").Value = 2639.879
$ws.Range("K138").Value = 5378.76
$ws.Range("L138").Value = 7919.637
$ws.Range("M138").Value = -238.7600000000002
$ws.Range("N138").Value = -18199.637

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1550
$ws.Range("I2").Value = 1500
$ws.Range("J2").Value = 1600
$ws.Range("K2").Value = 1500
$ws.Range("L2").Value = 1600
$ws.Range("M2").Value = -1387
$ws.Range("N2").Value = -1826
$ws.Range("H32").Value = 695248.9
$ws.Range("I32").Value = 799080.7
$ws.Range("K32").Value = 799080.7
$ws.Range("M32").Value = -798793.7
$ws.Range("H45").Value = 4204.5625
$ws.Range("J45").Value = 4132.6
$ws.Range("L45").Value = 4132.6
$ws.Range("N45").Value = -4886.6
$ws.Range("H74").Value = 850.25
$ws.Range("I74").Value = 710.2759
$ws.Range("J74").Value = 1430.1428
$ws.Range("K74").Value = 710.2759
$ws.Range("L74").Value = 1430.1428
$ws.Range("M74").Value = 163.7241
$ws.Range("N74").Value = -3178.1428
$ws.Range("H77").Value = 850.25
$ws.Range("I77").Value = 710.2759
$ws.Range("J77").Value = 1430.1428
$ws.Range("K77").Value = 3551.3795
$ws.Range("L77").Value = 7150.714
$ws.Range("M77").Value = 816.6205
$ws.Range("N77").Value = -15886.714
$ws.Range("H102").Value = 3442.2354
$ws.Range("I102").Value = 3647.8667
$ws.Range("K102").Value = 3647.8667
$ws.Range("M102").Value = -2025.8667
$ws.Range("H110").Value = 1982.7142
$ws.Range("I110").Value = 2098.2
$ws.Range("J110").Value = 1694
$ws.Range("K110").Value = 2098.2
$ws.Range("L110").Value = 1694
$ws.Range("M110").Value = -53.19999999999982
$ws.Range("N110").Value = -5784
$ws.Range("H116").Value = 1550
$ws.Range("I116").Value = 1500
$ws.Range("J116").Value = 1600
$ws.Range("K116").Value = 1500
$ws.Range("L116").Value = 1600
$ws.Range("M116").Value = 794
$ws.Range("N116").Value = -6188
$ws.Range("H132").Value = 2927.386
$ws.Range("I132").Value = 1996.3715
$ws.Range("J132").Value = 4408.5454
$ws.Range("K132").Value = 5989.1145
$ws.Range("L132").Value = 13225.6362
$ws.Range("M132").Value = -3459.1145
$ws.Range("N132").Value = -18285.6362

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1550
$ws.Range("I3").Value = 1500
$ws.Range("J3").Value = 1600
$ws.Range("K3").Value = 1500
$ws.Range("L3").Value = 1600
$ws.Range("M3").Value = -1386
$ws.Range("N3").Value = -1828
$ws.Range("H131").Value = 40136.668
$ws.Range("J131").Value = 40136.668
$ws.Range("L131").Value = 40136.668
$ws.Range("N131").Value = -50216.668
$ws.Range("H133").Value = 38718.32
$ws.Range("J133").Value = 38718.32
$ws.Range("L133").Value = 38718.32
$ws.Range("N133").Value = -48838.32

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 797.25
$ws.Range("I16").Value = 797.25
$ws.Range("K16").Value = 797.25
$ws.Range("M16").Value = -510.25
$ws.Range("H105").Value = 3000
$ws.Range("I105").Value = 0
$ws.Range("J105").Value = 3000
$ws.Range("K105").Value = 0
$ws.Range("L105").ClearContents()
$ws.Range("M105").Value = 3000
$ws.Range("N105").Value = -6494
$ws.Range("H113").Value = 797.25
$ws.Range("I113").Value = 797.25
$ws.Range("K113").Value = 797.25
$ws.Range("M113").Value = 1372.75
$ws.Range("H121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("L121").ClearContents()
$ws.Range("N121").Value = 0
$ws.Range("H122").Value = 1871.1951
$ws.Range("I122").Value = 1882.9412
$ws.Range("K122").Value = 5648.8236
$ws.Range("M122").Value = -3198.8236
$ws.Range("H132").Value = 7938738
$ws.Range("I132").Value = 1894.4546
$ws.Range("J132").Value = 16669266
$ws.Range("K132").Value = 5683.3638
$ws.Range("L132").Value = 50007798
$ws.Range("M132").Value = -3153.3638
$ws.Range("N132").Value = -50012858

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 934.05554
$ws.Range("J5").Value = 1966.8
$ws.Range("L5").Value = 5900.4
$ws.Range("N5").Value = -6124.4
$ws.Range("H20").Value = 1000
$ws.Range("I20").Value = 1000
$ws.Range("J20").Value = 1000
$ws.Range("K20").Value = 3000
$ws.Range("L20").Value = 3000
$ws.Range("M20").Value = -2773
$ws.Range("N20").Value = -3454
$ws.Range("H113").Value = 831.32556
$ws.Range("I113").Value = 609.6111
$ws.Range("J113").Value = 990.96
$ws.Range("K113").Value = 1828.8333
$ws.Range("L113").Value = 2972.88
$ws.Range("M113").Value = 341.1667000000002
$ws.Range("N113").Value = -7312.88
$ws.Range("H134").Value = 6870.0645
$ws.Range("I134").Value = 4329.5713
$ws.Range("J134").Value = 7611.0415
$ws.Range("K134").Value = 12988.7139
$ws.Range("L134").Value = 22833.1245
$ws.Range("M134").Value = -7918.713899999999
$ws.Range("N134").Value = -32973.12450000001
$ws.Range("H135").Value = 934.05554
$ws.Range("J135").Value = 1966.8
$ws.Range("L135").Value = 17701.2
$ws.Range("N135").Value = -22771.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 3997.5
$ws.Range("I46").Value = 1980.3334
$ws.Range("J46").Value = 4285.6665
$ws.Range("K46").Value = 1980.3334
$ws.Range("L46").Value = 4285.6665
$ws.Range("M46").Value = -1824.3334
$ws.Range("N46").Value = -4597.6665
$ws.Range("H97").Value = 1971
$ws.Range("I97").Value = 1735
$ws.Range("J97").Value = 2325
$ws.Range("K97").Value = 1735
$ws.Range("L97").Value = 2325
$ws.Range("M97").Value = -1239
$ws.Range("N97").Value = -3317
$ws.Range("H107").Value = 598.8
$ws.Range("I107").Value = 498
$ws.Range("J107").Value = 750
$ws.Range("K107").Value = 498
$ws.Range("L107").Value = 750
$ws.Range("M107").Value = 1422
$ws.Range("N107").Value = -4590
$ws.Range("H113").Value = 168035.33
$ws.Range("I113").Value = 250553
$ws.Range("K113").Value = 250553
$ws.Range("M113").Value = -248383

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 145557.72
$ws.Range("I40").Value = 202380.8
$ws.Range("J40").Value = 3500
$ws.Range("K40").Value = 202380.8
$ws.Range("L40").Value = 3500
$ws.Range("M40").Value = -202244.8
$ws.Range("N40").Value = -3772
$ws.Range("H93").Value = 5828.524
$ws.Range("I93").Value = 6967.353
$ws.Range("J93").Value = 988.5
$ws.Range("K93").Value = 6967.353
$ws.Range("L93").Value = 988.5
$ws.Range("M93").Value = -5719.353
$ws.Range("N93").Value = -3484.5
$ws.Range("H132").Value = 2215.0278
$ws.Range("I132").Value = 1483.2916
$ws.Range("J132").Value = 3678.5
$ws.Range("K132").Value = 4449.8748
$ws.Range("L132").Value = 11035.5
$ws.Range("M132").Value = -1919.8748
$ws.Range("N132").Value = -16095.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 4506213.5
$ws.Range("I132").Value = 1552
$ws.Range("J132").Value = 8773787
$ws.Range("K132").Value = 4656
$ws.Range("L132").Value = 26321361
$ws.Range("M132").Value = -2126
$ws.Range("N132").Value = -26326421
